$wb = $excel.ActiveWorkbook

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1830.375
$ws.Range("I101").Value = 243.25
$ws.Range("J101").Value = 3417.5
$ws.Range("K101").Value = 729.75
$ws.Range("L101").Value = 10252.5
$ws.Range("M101").Value = 892.25
$ws.Range("N101").Value = -13496.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4253.43
$ws.Range("I32").Value = 3961.361
$ws.Range("J32").Value = 5254.8096
$ws.Range("K32").Value = 3961.361
$ws.Range("L32").Value = 5254.8096
$ws.Range("M32").Value = -3674.361
$ws.Range("N32").Value = -5828.8096

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 50001836
$ws.Range("I63").Value = 50001836
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 50001836
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -50001150

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 50001836
$ws.Range("I66").Value = 50001836
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 250009180
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -250005748

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2967.3157
$ws.Range("I132").Value = 1750.1305
$ws.Range("J132").Value = 4833.6665
$ws.Range("K132").Value = 5250.3915
$ws.Range("L132").Value = 14500.9995
$ws.Range("M132").Value = -2720.3915
$ws.Range("N132").Value = -19560.9995

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27399.818
$ws.Range("I82").Value = 19499.5
$ws.Range("J82").Value = 29155.445
$ws.Range("K82").Value = 19499.5
$ws.Range("L82").Value = 29155.445
$ws.Range("M82").Value = -19116.5
$ws.Range("N82").Value = -29921.445

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 27399.818
$ws.Range("I85").Value = 19499.5
$ws.Range("J85").Value = 29155.445
$ws.Range("K85").Value = 19499.5
$ws.Range("L85").Value = 29155.445
$ws.Range("M85").Value = -18173.5
$ws.Range("N85").Value = -31807.445

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2227.75
$ws.Range("I86").Value = 1963.3
$ws.Range("J86").Value = 3550
$ws.Range("K86").Value = 1963.3
$ws.Range("L86").Value = 3550
$ws.Range("M86").Value = -840.3
$ws.Range("N86").Value = -5796

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2227.75
$ws.Range("I89").Value = 1963.3
$ws.Range("J89").Value = 3550
$ws.Range("K89").Value = 9816.5
$ws.Range("L89").Value = 17750
$ws.Range("M89").Value = -4200.5
$ws.Range("N89").Value = -28982

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10420702
$ws.Range("I99").Value = 1712.6666
$ws.Range("K99").Value = 1712.6666
$ws.Range("M99").Value = -214.6666

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10420702
$ws.Range("I126").Value = 1712.6666
$ws.Range("K126").Value = 5137.9998
$ws.Range("M126").Value = -2667.9998

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1292.25
$ws.Range("I68").Value = 889.61536
$ws.Range("J68").Value = 1728.4375
$ws.Range("K68").Value = 2668.84608
$ws.Range("L68").Value = 5185.3125
$ws.Range("M68").Value = -1857.84608
$ws.Range("N68").Value = -6807.3125

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1292.25
$ws.Range("I71").Value = 889.61536
$ws.Range("J71").Value = 1728.4375
$ws.Range("K71").Value = 8006.53824
$ws.Range("L71").Value = 15555.9375
$ws.Range("M71").Value = -3950.53824
$ws.Range("N71").Value = -23667.9375

# CUL row 112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2431.6
$ws.Range("J112").Value = 3357.1428
$ws.Range("L112").Value = 10071.4284
$ws.Range("N112").Value = -12287.4284

# CUL row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 298000
$ws.Range("I128").Value = 298000
$ws.Range("K128").Value = 894000
$ws.Range("M128").Value = -889020

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13757522
$ws.Range("J131").Value = 15386074
$ws.Range("L131").Value = 46158222
$ws.Range("N131").Value = -46168302

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 7314.5
$ws.Range("I140").Value = 7629.952
$ws.Range("J140").Value = 690
$ws.Range("K140").Value = 22889.856
$ws.Range("L140").Value = 2070
$ws.Range("M140").Value = -17709.856
$ws.Range("N140").Value = -12430

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 884.1429000000001
$ws.Range("I97").Value = 884.1429000000001
$ws.Range("K97").Value = 884.1429000000001
$ws.Range("M97").Value = -388.1429000000001

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3515.3809
$ws.Range("I132").Value = 8650
$ws.Range("J132").Value = 2974.8948
$ws.Range("K132").Value = 25950
$ws.Range("L132").Value = 8924.6844
$ws.Range("M132").Value = -23420
$ws.Range("N132").Value = -13984.6844

# GSM row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 78520
$ws.Range("J137").Value = 78520
$ws.Range("L137").Value = 78520
$ws.Range("N137").Value = -88720

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 169867.33
$ws.Range("I7").Value = 502502
$ws.Range("J7").Value = 3550
$ws.Range("K7").Value = 502502
$ws.Range("L7").Value = 3550
$ws.Range("M7").Value = -502390
$ws.Range("N7").Value = -3774

# LTW row 17
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 788.3333
$ws.Range("I17").Value = 788.3333
$ws.Range("K17").Value = 788.3333
$ws.Range("M17").Value = -618.3333

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2121.7144
$ws.Range("I61").Value = 1808.25
$ws.Range("J61").Value = 4002.5
$ws.Range("K61").Value = 1808.25
$ws.Range("L61").Value = 4002.5
$ws.Range("M61").Value = -1606.25
$ws.Range("N61").Value = -4406.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2121.7144
$ws.Range("I113").Value = 1808.25
$ws.Range("J113").Value = 4002.5
$ws.Range("K113").Value = 1808.25
$ws.Range("L113").Value = 4002.5
$ws.Range("M113").Value = 361.75
$ws.Range("N113").Value = -8342.5

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 169867.33
$ws.Range("I126").Value = 502502
$ws.Range("J126").Value = 3550
$ws.Range("K126").Value = 1507506
$ws.Range("L126").Value = 10650
$ws.Range("M126").Value = -1505036
$ws.Range("N126").Value = -15590

# WVR row 17
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 70005
$ws.Range("J17").Value = 70005
$ws.Range("L17").Value = 70005
$ws.Range("N17").Value = -70349

# WVR row 51
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 669.4865
$ws.Range("I126").Value = 539.73334
$ws.Range("J126").Value = 1225.5714
$ws.Range("K126").Value = 1619.20002
$ws.Range("L126").Value = 3676.7142
$ws.Range("M126").Value = 850.79998
$ws.Range("N126").Value = -8616.7142
